# Fixed naive component forecaster bug - Presentation state 11.02.
#
# A new "1-quarter-ahead" naive forecast error value is now available for
# every vintage date (row). It needs to be inserted as the new column B
# ("Q0") for each row, pushing the previously computed errors (old Q0..Q8)
# one column to the right (becoming the new Q1..Q9), with whatever value
# used to sit in the last populated column of the row (which would now
# overflow past column K / Q9) being dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly computed "Q0" (1-quarter-ahead) naive forecast error per row (2-20).
$newValues = @{
    2  = 0.7916129955631771
    3  = -3.727363316492332
    4  = 0.376932102669816
    5  = 1.207578635508109
    6  = -0.9264868865757077
    7  = 0.3770345820039356
    8  = -0.4275923834192769
    9  = 0.324932645901923
    10 = -0.04071760298358112
    11 = 0.3721869518844864
    12 = -0.1524291232873974
    13 = -1.030518528898312
    14 = 0.4742145784871607
    15 = 0.3556547466179877
    16 = 0.3126006297022321
    17 = 0.3812981176718321
    18 = -0.716162849403934
    19 = 0.506656010950813
    20 = -0.343237405067616
}

$firstDataCol = 2   # column B
$lastDataCol  = 11  # column K

for ($row = 2; $row -le 20; $row++) {

    # Determine the last populated column (B..K) in this row before editing it.
    $lastCol = $firstDataCol - 1
    for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
        if ($ws.Cells.Item($row, $c).Value2 -ne $null) {
            $lastCol = $c
        }
    }

    # Snapshot the existing values (B..lastCol) before they get overwritten.
    $oldVals = @()
    for ($c = $firstDataCol; $c -le $lastCol; $c++) {
        $oldVals += $ws.Cells.Item($row, $c).Value2
    }

    # Insert the new value into column B.
    $ws.Cells.Item($row, $firstDataCol).Value2 = $newValues[$row]

    # Shift the previously-existing values one column to the right,
    # dropping the final value if it would spill past column K.
    for ($i = 0; $i -lt $oldVals.Count; $i++) {
        $destCol = $firstDataCol + 1 + $i
        if ($destCol -le $lastDataCol) {
            $ws.Cells.Item($row, $destCol).Value2 = $oldVals[$i]
        }
    }
}
